# Update sheet/workbook + vendor name, and sales rows for Campanha 67.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Terezinha França"

$ws.Cells.Item(3,1).Value = "Terezinha De Jesus Jordão De França"
$ws.Cells.Item(3,2).Value = 46321
$ws.Cells.Item(3,3).Value = "Ótica Daily"
$ws.Cells.Item(3,4).Value = "Ervio Jose Boconcelo"
$ws.Cells.Item(3,5).Value = "21/12/2024"
$ws.Cells.Item(3,6).Value = "R$ 1.200,00"
$ws.Cells.Item(3,7).Value = "Convertido"

$ws.Cells.Item(4,1).Value = "Terezinha De Jesus Jordão De França"
$ws.Cells.Item(4,2).Value = 46374
$ws.Cells.Item(4,3).Value = "Ótica Daily"
$ws.Cells.Item(4,4).Value = "Dejanira De Alcantara Pereira"
$ws.Cells.Item(4,5).Value = "08/01/2025"
$ws.Cells.Item(4,6).Value = "R$ 1.000,00"
$ws.Cells.Item(4,7).Value = "Atrasado"

$ws.Cells.Item(5,1).Value = "Terezinha De Jesus Jordão De França"
$ws.Cells.Item(5,2).Value = 46450
$ws.Cells.Item(5,3).Value = "Ótica Daily"
$ws.Cells.Item(5,4).Value = "Niva Pancotti Mendonca"
$ws.Cells.Item(5,5).Value = "11/01/2025"
$ws.Cells.Item(5,6).Value = "R$ 470,00"
$ws.Cells.Item(5,7).Value = "Atrasado"

$ws.Cells.Item(6,1).Value = "Terezinha De Jesus Jordão De França"
$ws.Cells.Item(6,2).Value = 46481
$ws.Cells.Item(6,3).Value = "Ótica Daily"
$ws.Cells.Item(6,4).Value = "Michel Simion Mlechecov"
$ws.Cells.Item(6,5).Value = "14/01/2025"
$ws.Cells.Item(6,6).Value = "R$ 290,00"
$ws.Cells.Item(6,7).Value = "Convertido"

$ws.Cells.Item(7,1).Value = "Terezinha De Jesus Jordão De França"
$ws.Cells.Item(7,2).Value = 46525
$ws.Cells.Item(7,3).Value = "Ótica Daily"
$ws.Cells.Item(7,4).Value = "Alessandra Simone Silva Fantucci"
$ws.Cells.Item(7,5).Value = "17/01/2025"
$ws.Cells.Item(7,6).Value = "R$ 180,00"
$ws.Cells.Item(7,7).Value = "Atrasado"

$ws.Cells.Item(8,6).Value = "Total Pendente:"
$ws.Cells.Item(8,7).Value = "R$ 1.650,00"
$ws.Cells.Item(9,6).Value = "Total Convertido:"
$ws.Cells.Item(9,7).Value = "R$ 1.490,00"
$ws.Cells.Item(10,6).Value = "Total:"
$ws.Cells.Item(10,7).Value = "R$ 3.140,00"
